$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Change 1: Pos-Condicoes section - "Veiculo registrado no estacionamento" +
# "." (two runs) become a single run "Cancela liberada para entrada do
# veiculo. Veiculo registrado no estacionamento."
# ---------------------------------------------------------------------------
$found = $d.Content.Find.Execute(
    "Veículo registrado no estacionamento.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Cancela liberada para entrada do veículo. Veículo registrado no estacionamento.",
    2)

# ---------------------------------------------------------------------------
# Change 2: Fluxo Basico - insert two new "Titulo3" steps right before the
# closing "O caso de uso e encerrado." step:
#   - "Sistema informa que a entrada do veiculo foi registrada"
#   - "Cancela libera a entrada do veiculo. " (kept in the same paragraph as
#     the pre-existing "O caso de uso e encerrado." run)
# ---------------------------------------------------------------------------
$targetIdx = 0
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text -like "*O caso de uso é encerrado*") {
        $targetIdx = $i
        break
    }
}

$closingPara = $d.Paragraphs.Item($targetIdx)

# Split the closing paragraph: insert "Sistema informa..." + a paragraph
# break right before the existing (untouched) "O caso de uso..." run. This
# leaves that original run alone in its own paragraph, one position later.
$splitRange = $closingPara.Range.Duplicate
$splitRange.Collapse(1)
$splitRange.Text = "Sistema informa que a entrada do veículo foi registrada`r"

# The paragraph that now holds only the original "O caso de uso..." run.
$closingPara = $d.Paragraphs.Item($targetIdx + 1)

# Add an empty paragraph right before it (inherits the same Titulo3 style)
# so we can type the new lead-in sentence without touching/merging the
# pre-existing run.
$closingPara.Range.InsertParagraphBefore()
$leadInPara = $d.Paragraphs.Item($targetIdx + 1)
$closingPara = $d.Paragraphs.Item($targetIdx + 2)

$sel = $word.Selection
$sel.SetRange($leadInPara.Range.End - 1, $leadInPara.Range.End - 1)
$sel.TypeText("Cancela libera a entrada do veículo. ")

# Join the lead-in paragraph back with the closing paragraph (removes the
# paragraph mark between them) so both runs end up in the same paragraph,
# exactly as in the target document, without merging the two runs.
$markRange = $d.Range($leadInPara.Range.End - 1, $leadInPara.Range.End)
$markRange.Delete()
